$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. Update the existing registration row's e-mail address (D2) and keep the
#    hyperlink pointing at the same mailbox, preserving the bordered
#    hyperlink cell look that was already there.
# ---------------------------------------------------------------------------
$ws.Range("D2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial($xlPasteFormats) | Out-Null   # stash original look

$ws.Range("D2").Hyperlinks.Delete()
$ws.Range("D2").Value = "Chaitanyatawade1@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Chaitanyatawade1@gmail.com") | Out-Null

$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D2").PasteSpecial($xlPasteFormats) | Out-Null   # restore bordered hyperlink look
$ws.Range("Z1").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 2. Append five more registration rows (3-7), each with its own e-mail
#    hyperlink. Some keep the bordered look (copied from D2), others use the
#    plain "Hyperlink" cell style without borders.
# ---------------------------------------------------------------------------
function Add-Registration($row, $gender, $first, $last, $email, $bordered) {
    $ws.Cells.Item($row, 1).Value = $gender
    $ws.Cells.Item($row, 2).Value = $first
    $ws.Cells.Item($row, 3).Value = $last
    $ws.Cells.Item($row, 5).Value = "Abc@123"
    $ws.Cells.Item($row, 6).Value = "Abc@123"

    $d = $ws.Cells.Item($row, 4)
    $d.Value = $email
    $ws.Hyperlinks.Add($d, "mailto:" + $email) | Out-Null

    if ($bordered) {
        $ws.Range("D2").Copy() | Out-Null
        $d.PasteSpecial($xlPasteFormats) | Out-Null
    } else {
        $d.Style = "Hyperlink"
    }
}

Add-Registration 3 "Male"   "Om"     "Tawade" "Chaitanyatawade2@gmail.com" $false
Add-Registration 4 "Male"   "Nitesh" "Tawade" "Chaitanyatawade1@gmail.com" $true
Add-Registration 6 "Male"   "Deepak" "Tawade" "Chaitanyatawade1@gmail.com" $true
Add-Registration 5 "Female" "Jyoti"  "Tawade" "Chaitanyatawade2@gmail.com" $false
Add-Registration 7 "Male"   "TEST"   "Tawade" "Chaitanyatawade2@gmail.com" $false

# ---------------------------------------------------------------------------
# 3. Borders on the plain (non-hyperlink) cells of the new rows - copy the
#    existing bordered look from row 2 rather than building a fresh border
#    (keeps the same style / colour as the rest of the sheet).
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:C7").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3:F7").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("A5").Select() | Out-Null
